$d = $word.ActiveDocument

# Helper: build a brand-new Range object from a Paragraph's current
# Start/End. Reusing a Range obtained straight off Paragraph.Next() as the
# receiver of a later Find.Execute call can make the search spill outside
# the paragraph even with Wrap=wdFindStop, so we always re-wrap the
# boundaries in a fresh Document.Range(...) before searching in it.
function Get-FreshRange($para) {
    return $d.Range($para.Range.Start, $para.Range.End)
}

# --- 1. SRS_Home_2 requirement body --------------------------------------
# "The user can navigate to the booking and rating system by clicking on
# any trip." -> "The registered user can navigate to the booking and
# rating system by clicking on any trip."
# Anchor on the "SRS_Home_2:" label (unique in the doc) and walk forward
# two paragraphs (label -> blank line -> body) so we only touch this one
# occurrence of the very common phrase "The us...".
$lbl1 = $d.Content
$foundLbl1 = $lbl1.Find.Execute("SRS_Home_2:")
$ok1 = $false
if ($foundLbl1) {
    $labelPara1 = $lbl1.Paragraphs(1)
    $bodyPara1 = $labelPara1.Next().Next()
    $bodyRange1 = Get-FreshRange $bodyPara1
    $ok1 = $bodyRange1.Find.Execute("The us", $true, $false, $false, $false, $false, `
        $true, 0, $false, "The registered us", 2)
}

# --- 2. SRS_Admn_6 requirement body ---------------------------------------
# "The program should display a message containing "User added
# Successfully"" -> "When the admin deletes user, the program should
# display a message containing "User added Successfully""
# Two paragraphs in the doc start with "The program should display a
# message containing", so anchor on the "SRS_Admn_6:" label immediately
# before the one we want.
$lbl2 = $d.Content
$foundLbl2 = $lbl2.Find.Execute(" SRS_Admn_6:")
$ok2 = $false
if ($foundLbl2) {
    $labelPara2 = $lbl2.Paragraphs(1)
    $bodyPara2 = $labelPara2.Next()
    $bodyRange2 = Get-FreshRange $bodyPara2
    $ok2 = $bodyRange2.Find.Execute("The program", $true, $false, $false, $false, $false, `
        $true, 0, $false, "When the admin deletes user, the program", 2)
}

# --- 3. Remove the SRS_Admn_12 requirement entirely (label + body) -------
$lbl3 = $d.Content
$foundLbl3 = $lbl3.Find.Execute("SRS_Admn_12:")
if ($foundLbl3) {
    $labelPara3 = $lbl3.Paragraphs(1)
    $bodyPara3 = $labelPara3.Next()
    $delRange = $d.Range($labelPara3.Range.Start, $bodyPara3.Range.End)
    $delRange.Delete()
}

# --- 4. SRS_BK_1 requirement body -----------------------------------------
# "The user can book any flight and a message appears containing "Flight
# Booked Successfully"." -> "Registered user can book any flight and a
# message appears containing "Flight Booked Successfully"."
# "The user can book any flight and" is unique document-wide, so a plain
# whole-document search is safe here.
$rng4 = $d.Content
$ok4 = $rng4.Find.Execute("The user can book", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Registered user can book", 2)

Write-Output "foundLbl1=$foundLbl1 ok1=$ok1 foundLbl2=$foundLbl2 ok2=$ok2 foundLbl3=$foundLbl3 ok4=$ok4"
